# This workbook is a point-in-time stock report; the update replaces it with a
# newer extract. Net effect on the data rows falls into three patterns:
#   1) A pair of adjacent rows for the same item (same code/description/rate) has
#      its stock-code / rate / qty / value columns swapped (re-ordered receipts).
#   2) A single row keeps its code/rate but its qty (and qty*rate value) changes.
#   3) Each section/grand "Sub Total:" (and the final "Grand Total:") row is
#      recomputed to the new sum of its section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Adjacent same-item rows: swap Code(B) / Rate(E) / Qty(F) / Value(G) ---
# Rows 33 & 34
$ws.Range("B33").Value = 66452
$ws.Range("F33").Value = 70
$ws.Range("G33").Value = 2151.8
$ws.Range("B34").Value = 51755
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 30.74

# Rows 151 & 152
$ws.Range("B151").Value = 64196
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 32143.58
$ws.Range("B152").Value = 65258
$ws.Range("F152").Value = 2
$ws.Range("G152").Value = 64287.16

# Rows 388 & 389
$ws.Range("B388").Value = 57802
$ws.Range("E388").Value = 162.71
$ws.Range("F388").Value = -79
$ws.Range("G388").Value = -11334.92
$ws.Range("B389").Value = 63531
$ws.Range("E389").Value = 152.53
$ws.Range("F389").Value = 39
$ws.Range("G389").Value = 5595.72

# Rows 400 & 401
$ws.Range("B400").Value = 60325
$ws.Range("E400").Value = 151.57
$ws.Range("F400").Value = -102
$ws.Range("G400").Value = -12939.72
$ws.Range("B401").Value = 63560
$ws.Range("E401").Value = 134.87
$ws.Range("F401").Value = 1
$ws.Range("G401").Value = 126.86

# Rows 553 & 554
$ws.Range("B553").Value = 65066
$ws.Range("E553").Value = 13.61
$ws.Range("F553").Value = 90
$ws.Range("G553").Value = 1152.9
$ws.Range("B554").Value = 53263
$ws.Range("E554").Value = 15.29
$ws.Range("F554").Value = -309
$ws.Range("G554").Value = -3958.29

# Rows 562 & 563
$ws.Range("B562").Value = 45718
$ws.Range("E562").Value = 19.38
$ws.Range("F562").Value = -294
$ws.Range("G562").Value = -4768.68
$ws.Range("B563").Value = 64927
$ws.Range("E563").Value = 17.26
$ws.Range("F563").Value = 106
$ws.Range("G563").Value = 1719.32

# Rows 567 & 568
$ws.Range("B567").Value = 64925
$ws.Range("E567").Value = 13.97
$ws.Range("F567").Value = 111
$ws.Range("G567").Value = 1459.65
$ws.Range("B568").Value = 45709
$ws.Range("E568").Value = 15.69
$ws.Range("F568").Value = -300
$ws.Range("G568").Value = -3945

# Rows 572 & 573
$ws.Range("B572").Value = 53595
$ws.Range("E572").Value = 17.61
$ws.Range("F572").Value = -335
$ws.Range("G572").Value = -4934.55
$ws.Range("B573").Value = 65067
$ws.Range("E573").Value = 15.65
$ws.Range("F573").Value = 126
$ws.Range("G573").Value = 1855.98

# Rows 672 & 673
$ws.Range("B672").Value = 60022
$ws.Range("E672").Value = 37.22
$ws.Range("F672").Value = -113
$ws.Range("G672").Value = -3709.79
$ws.Range("B673").Value = 64830
$ws.Range("E673").Value = 34.9
$ws.Range("F673").Value = 91
$ws.Range("G673").Value = 2987.53

# --- 2) Quantity (F) and recomputed Value (G) updates on individual rows ---
$ws.Range("F45").Value = 585
$ws.Range("G45").Value = 21299.85
$ws.Range("F48").Value = 276
$ws.Range("G48").Value = 53237.64
$ws.Range("F59").Value = 225
$ws.Range("G59").Value = 21046.5
$ws.Range("F60").Value = 165
$ws.Range("G60").Value = 9735
$ws.Range("F73").Value = 40
$ws.Range("G73").Value = 1338.4
$ws.Range("F114").Value = 24
$ws.Range("G114").Value = 8923.200000000001
$ws.Range("F124").Value = 87
$ws.Range("G124").Value = 7061.79
$ws.Range("F177").Value = 112
$ws.Range("G177").Value = 5890.08
$ws.Range("F239").Value = 70
$ws.Range("G239").Value = 3029.6
$ws.Range("F243").Value = 33
$ws.Range("G243").Value = 478.83
$ws.Range("F285").Value = 83
$ws.Range("G285").Value = 5243.94
$ws.Range("F287").Value = 2325
$ws.Range("G287").Value = 43012.5
$ws.Range("F324").Value = 11
$ws.Range("G324").Value = 5019.19
$ws.Range("F339").Value = 218
$ws.Range("G339").Value = 9260.639999999999
$ws.Range("F342").Value = 19
$ws.Range("G342").Value = 1628.68
$ws.Range("F362").Value = 239
$ws.Range("G362").Value = 11201.93
$ws.Range("F409").Value = 137
$ws.Range("G409").Value = 20710.29
$ws.Range("F540").Value = 0
$ws.Range("G540").Value = 0
$ws.Range("F543").Value = 0
$ws.Range("G543").Value = 0
$ws.Range("F578").Value = 18
$ws.Range("G578").Value = 901.4400000000001
$ws.Range("F581").Value = 158
$ws.Range("G581").Value = 7648.78
$ws.Range("F629").Value = 45
$ws.Range("G629").Value = 2515.05
$ws.Range("F635").Value = 19
$ws.Range("G635").Value = 2951.08
$ws.Range("F680").Value = 383
$ws.Range("G680").Value = 38250.21
$ws.Range("F684").Value = 34
$ws.Range("G684").Value = 9719.58
$ws.Range("F700").Value = 115
$ws.Range("G700").Value = 3807.65
$ws.Range("F704").Value = 274
$ws.Range("G704").Value = 9072.139999999999
$ws.Range("F706").Value = 67
$ws.Range("G706").Value = 2015.36
$ws.Range("F798").Value = 0
$ws.Range("G798").Value = 0
$ws.Range("F799").Value = 264
$ws.Range("G799").Value = 35138.4
$ws.Range("F805").Value = 27
$ws.Range("G805").Value = 4425.03
$ws.Range("F806").Value = 4
$ws.Range("G806").Value = 435.24
$ws.Range("F807").Value = 167
$ws.Range("G807").Value = 18171.27
$ws.Range("F812").Value = 46
$ws.Range("G812").Value = 6731.18
$ws.Range("F815").Value = 105
$ws.Range("G815").Value = 15233.4
$ws.Range("F818").Value = 10
$ws.Range("G818").Value = 1559.2
$ws.Range("F819").Value = 113
$ws.Range("G819").Value = 5437.56
$ws.Range("F825").Value = 459
$ws.Range("G825").Value = 36072.81
$ws.Range("F832").Value = 99
$ws.Range("G832").Value = 4673.79
$ws.Range("F842").Value = 5
$ws.Range("G842").Value = 33688.3
$ws.Range("F844").Value = 5
$ws.Range("G844").Value = 70780.8
$ws.Range("F891").Value = 1866
$ws.Range("G891").Value = 304363.26
$ws.Range("F933").Value = 26
$ws.Range("G933").Value = 5565.56

# --- 3) Section Sub Total / Grand Total (B) recomputed from the new row values ---
$ws.Range("B74").Value = 291726.41
$ws.Range("B116").Value = 31714.99
$ws.Range("B147").Value = 99962.03
$ws.Range("B184").Value = 30666.8
$ws.Range("B250").Value = 102935.89
$ws.Range("B294").Value = 65414.49
$ws.Range("B375").Value = 181303.67
$ws.Range("B412").Value = 55530.3
$ws.Range("B548").Value = 21317
$ws.Range("B587").Value = 43196.88
$ws.Range("B636").Value = 72200.3
$ws.Range("B695").Value = 190336.07
$ws.Range("B708").Value = 42247.7
$ws.Range("B802").Value = 37120.65
$ws.Range("B838").Value = 341778.97
$ws.Range("B845").Value = 326779.8
$ws.Range("B897").Value = 353551.87
$ws.Range("B936").Value = 123556.58
$ws.Range("B942").Value = 5268470.15
$ws.Range("B943").Value = 5268470.15
